$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Overview" sheet: b.md.md row (row 3) flips from "handed back" to
#    "ready for handoff" in both the zh-cn and de-de status columns.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# Helper data per language sheet: new handoff file name + new handoff
# datetime for the b.md.md row (row 3).
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandoffFile = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"; HandoffTime = "2016-01-25 07:01:08" },
    @{ Name = "de-de"; HandoffFile = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"; HandoffTime = "2016-01-25 07:01:20" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (B) for the b.md.md row goes from the handed-back text
    # to "Ready for handoff".
    $ws.Range("B3").Value = "Ready for handoff"

    # Latest Handoff File (C) / Latest Handoff Datetime (D) now reference the
    # brand new handoff package for b.md.md.
    $ws.Range("C3").Value = $lang.HandoffFile
    $ws.Range("D3").Value = $lang.HandoffTime

    # The hyperlink sitting on C3 needs its visible text updated to match -
    # the underlying target URL/relationship id is unchanged, so update the
    # existing Hyperlink object in place instead of removing/re-adding it.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$3') {
            $h.TextToDisplay = $lang.HandoffFile
        }
    }
}
